$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay as plain text, matching source formatting,
# even when the new value looks numeric (e.g. "535.26" or "0.0000326").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.934.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.026.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.023.22"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.697"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.95"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000326"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.665.94"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.028.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.18"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.83"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.925.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.91"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "98.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.60"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +29.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.39"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.80"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.06"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +20.85%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.11"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +17.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.59"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "676.79"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "68.13"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.460"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.69%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.65%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.151"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.39"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.861.10"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.26%  "

# Row swaps: rows 39/40 (WEMIXToken <-> PEPE) and 43/44 (Dai <-> THORChain)
# swap Coin/Link/Price/Volume data between the row pairs, with updated Price/Volume values.
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0824"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.95%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.13"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +16.27%  "
$ws.Range("B44").Value = "Dai"
$ws.Range("C44").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.03%  "
